$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 41/42 swap content: Aptos now ranks above TheSandbox.
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.73%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7682"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.71%  "

# Price (D) and Volume(1h) (E) refreshes for every other row.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.293.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.56%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3966"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3936"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.388"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08579"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.317"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.022"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.48%  "

$ws.Range("E16").Value = "  +3.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.665.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.004"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.34%  "

$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("E23").Value = "  +2.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.296.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.539"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.105"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "142.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.360"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.002"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.534"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.848.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.058"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03069"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08267"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.862"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2761"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09289"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.440"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.93%  "

$ws.Range("E45").Value = "  +2.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.530"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.124"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("E48").Value = "  +0.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08406"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.262"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.39%  "

